$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the taxon-observation data between row 4 and row 6
# (same "Kratte masugn, Gstr" site, two different species records that
# traded places in the source export). Update only the cells that
# actually differ between the two rows, swapping row 4's values with
# row 6's values.

# --- Row 4 <- values that used to be on row 6 ---
$ws.Range("A4").Value = 111896689
$ws.Range("B4").Value = 90687
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 5964
$ws.Range("F4").Value = "Fjällig taggsvamp s.str."
$ws.Range("G4").Value = "Sarcodon imbricatus s.str."
$ws.Range("H4").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q4").Value = 575759.4008215864
$ws.Range("R4").Value = 6703742.042864766

# --- Row 6 <- values that used to be on row 4 ---
$ws.Range("A6").Value = 111896686
$ws.Range("B6").Value = 88966
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 5754
$ws.Range("F6").Value = "Gultoppig fingersvamp"
$ws.Range("G6").Value = "Ramaria testaceoflava"
$ws.Range("H6").Value = "(Bres.) Corner"
$ws.Range("Q6").Value = 575755.4490459578
$ws.Range("R6").Value = 6703741.960185026

# The (empty) "Bestämningsmetod" cell also moved from row 6 to row 4.
# Writing "" collapses a cell to fully blank in this engine, so use the
# classic text-quote-prefix trick to force an empty *text* cell, then
# strip the quote-prefix formatting it leaves behind.
$ws.Range("AF6").ClearContents()
$ws.Range("AF4").Value = "'"
$ws.Range("AF4").Style = "Normal"
